$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-register the new shared strings in the same order the target workbook
# introduces them (C22984, C128551, C31850 -> sharedStrings indices 9,10,11).
$ws.Range("M13").Value = "C22984"
$ws.Range("N14").Value = "C128551"
$ws.Range("M12").Value = "C31850"

# --- Row 12 : R2 (sel) input row (mirrors row 4, the R1 (sel) row) ---
$ws.Range("A4").Copy($ws.Range("A12"))
$ws.Range("A12").Value = "R2 (sel)"

$ws.Range("B4").Copy($ws.Range("B12"))
$ws.Range("B12").Value = 22000

$ws.Range("C4").Copy($ws.Range("C12"))
$ws.Range("C12").Value = "Ohm"

$ws.Range("M12").Value = "C31850"

# --- Row 13 : Vout (sel) = 1.2V row (mirrors row 5) ---
$ws.Range("A5").Copy($ws.Range("A13"))
$ws.Range("A13").Value = "Vout (sel)"

$ws.Range("B5").Copy($ws.Range("B13"))
$ws.Range("B13").Value = 1.2

$ws.Range("C5").Copy($ws.Range("C13"))
$ws.Range("C13").Value = "V"

$ws.Range("D5").Copy($ws.Range("D13"))
$ws.Range("D13").Value = "R1 (calc)"

$ws.Range("E5").Copy($ws.Range("E13"))
$ws.Range("E13").Formula = "=`$B`$12*(`$B13/`$B`$3-1)"

$ws.Range("F5").Copy($ws.Range("F13"))
$ws.Range("F13").Value = "Ohm"

$ws.Range("G5").Copy($ws.Range("G13"))
$ws.Range("G13").Value = "R1 (sel)"

$ws.Range("H5").Copy($ws.Range("H13"))
$ws.Range("H13").Value = 30000

$ws.Range("I5").Copy($ws.Range("I13"))
$ws.Range("I13").Value = "Ohm"

$ws.Range("J5").Copy($ws.Range("J13"))
$ws.Range("J13").Value = "Vout (calc)"

$ws.Range("K5").Copy($ws.Range("K13"))
$ws.Range("K13").Formula = "=(`$H13/`$B`$12+1)*`$B`$3"

$ws.Range("L5").Copy($ws.Range("L13"))
$ws.Range("L13").Value = "V"

$ws.Range("M13").Value = "C22984"

# --- Row 14 : Vout (sel) = 5.5V row (mirrors row 6, which has a formula-driven H cell) ---
$ws.Range("A6").Copy($ws.Range("A14"))
$ws.Range("A14").Value = "Vout (sel)"

$ws.Range("B6").Copy($ws.Range("B14"))
$ws.Range("B14").Value = 5.5

$ws.Range("C6").Copy($ws.Range("C14"))
$ws.Range("C14").Value = "V"

$ws.Range("D6").Copy($ws.Range("D14"))
$ws.Range("D14").Value = "R1 (calc)"

$ws.Range("E6").Copy($ws.Range("E14"))
$ws.Range("E14").Formula = "=`$B`$12*(`$B14/`$B`$3-1)"

$ws.Range("F6").Copy($ws.Range("F14"))
$ws.Range("F14").Value = "Ohm"

$ws.Range("G6").Copy($ws.Range("G14"))
$ws.Range("G14").Value = "R1 (sel)"

$ws.Range("H6").Copy($ws.Range("H14"))
$ws.Range("H14").Formula = "=200000+30000"

$ws.Range("I6").Copy($ws.Range("I14"))
$ws.Range("I14").Value = "Ohm"

$ws.Range("J6").Copy($ws.Range("J14"))
$ws.Range("J14").Value = "Vout (calc)"

$ws.Range("K6").Copy($ws.Range("K14"))
$ws.Range("K14").Formula = "=(`$H14/`$B`$12+1)*`$B`$3"

$ws.Range("L6").Copy($ws.Range("L14"))
$ws.Range("L14").Value = "V"

$ws.Range("M14").Value = "C22984"
$ws.Range("N14").Value = "C128551"
